$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Salary/5000 -> Teaching/1500 (Date in C2 is unchanged) ---
$ws.Range("A2").Value = "Teaching"
$ws.Range("B2").Value = 1500

# --- Row 3 (new): Money / 2000 / 2026-01-16 ---
$ws.Range("A3").Value = "Money"
$ws.Range("B3").Value = 2000
$ws.Range("C3").Value = 46038.291712962964

# --- Row 4 (new): Birthday / 2500 / 2026-01-15 ---
$ws.Range("A4").Value = "Birthday"
$ws.Range("B4").Value = 2500
$ws.Range("C4").Value = 46037.291712962964

# Copy the existing date formatting from C2 onto the two new date cells so
# they share the same (date) style/number format instead of plain numbers.
$ws.Range("C2").Copy()
$ws.Range("C3:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
